$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.74403644308763
$ws.Range("C2").Value = 10.68104836919495
$ws.Range("D2").Value = 12.84230271379829
$ws.Range("E2").Value = 12.05414600804594
$ws.Range("G2").Value = 3.841133873346237
$ws.Range("I2").Value = 44.45218186239924
$ws.Range("J2").Value = 7.479373422515214
$ws.Range("K2").Value = 22.36385840546215
$ws.Range("L2").Value = 14.82135193464963
$ws.Range("N2").Value = 28.11438811150361
$ws.Range("B3").Value = 24.6667897445041
$ws.Range("C3").Value = 10.57657537048508
$ws.Range("D3").Value = 12.84536146111778
$ws.Range("E3").Value = 12.07593591476398
$ws.Range("G3").Value = 3.844784390681683
$ws.Range("I3").Value = 44.43295455265356
$ws.Range("J3").Value = 7.482834944910886
$ws.Range("K3").Value = 22.31151922590653
$ws.Range("L3").Value = 14.83598901566326
$ws.Range("N3").Value = 28.13092495601708
$ws.Range("B4").Value = 24.62524050407935
$ws.Range("C4").Value = 10.51482489021473
$ws.Range("D4").Value = 12.84951445484534
$ws.Range("E4").Value = 12.09070721792382
$ws.Range("G4").Value = 3.847142694716771
$ws.Range("I4").Value = 44.42703674037769
$ws.Range("J4").Value = 7.485103127229858
$ws.Range("K4").Value = 22.28441142537943
$ws.Range("L4").Value = 14.84725983334297
$ws.Range("N4").Value = 28.14283343510605
$ws.Range("B5").Value = 24.60979919940553
$ws.Range("C5").Value = 10.49028851762277
$ws.Range("D5").Value = 12.85177919602448
$ws.Range("E5").Value = 12.09707709092429
$ws.Range("G5").Value = 3.848133219835138
$ws.Range("I5").Value = 44.42610684186638
$ws.Range("J5").Value = 7.486063426015513
$ws.Range("K5").Value = 22.27463617755893
$ws.Range("L5").Value = 14.85242699179487
$ws.Range("N5").Value = 28.14812742474333
$ws.Range("B6").Value = 24.60732550279425
$ws.Range("C6").Value = 10.48625286764705
$ws.Range("D6").Value = 12.85218982903806
$ws.Range("E6").Value = 12.09815598010455
$ws.Range("G6").Value = 3.848299480522443
$ws.Range("I6").Value = 44.42604194303419
$ws.Range("J6").Value = 7.486225059612002
$ws.Range("K6").Value = 22.27308997244253
$ws.Range("L6").Value = 14.8533196725619
$ws.Range("N6").Value = 28.14903313192843
$ws.Range("B7").Value = 24.62502620882679
$ws.Range("C7").Value = 10.51449141177335
$ws.Range("D7").Value = 12.84954268019433
$ws.Range("E7").Value = 12.0907917048347
$ws.Range("G7").Value = 3.847155933709488
$ws.Range("I7").Value = 44.42701819935366
$ws.Range("J7").Value = 7.485115932287687
$ws.Range("K7").Value = 22.28427443695291
$ws.Range("L7").Value = 14.84732719461792
$ws.Range("N7").Value = 28.14290304553629
$ws.Range("B8").Value = 24.71618875739388
$ws.Range("C8").Value = 10.64454403793175
$ws.Range("D8").Value = 12.84288546102447
$ws.Range("E8").Value = 12.06137043523955
$ws.Range("G8").Value = 3.842368380601601
$ws.Range("I8").Value = 44.44433108126255
$ws.Range("J8").Value = 7.480537372332309
$ws.Range("K8").Value = 22.34477265421191
$ws.Range("L8").Value = 14.82592485935354
$ws.Range("N8").Value = 28.11972577009174
$ws.Range("B9").Value = 24.94103178333231
$ws.Range("C9").Value = 10.91743718911387
$ws.Range("D9").Value = 12.84786016541182
$ws.Range("E9").Value = 12.0147059967463
$ws.Range("G9").Value = 3.833902368734528
$ws.Range("I9").Value = 44.52493463723537
$ws.Range("J9").Value = 7.472687838394056
$ws.Range("K9").Value = 22.5029544780815
$ws.Range("L9").Value = 14.80207233875954
$ws.Range("N9").Value = 28.08820285189454
$ws.Range("B10").Value = 25.13342161770335
$ws.Range("C10").Value = 11.12720293639298
$ws.Range("D10").Value = 12.8624679674272
$ws.Range("E10").Value = 11.98712462888712
$ws.Range("G10").Value = 3.828237751588463
$ws.Range("I10").Value = 44.6124862817914
$ws.Range("J10").Value = 7.46760359420013
$ws.Range("K10").Value = 22.64272823899595
$ws.Range("L10").Value = 14.79558612070811
$ws.Range("N10").Value = 28.07354098843669
$ws.Range("B11").Value = 25.22662837242319
$ws.Range("C11").Value = 11.22428473514539
$ws.Range("D11").Value = 12.87147908822844
$ws.Range("E11").Value = 11.97602789644083
$ws.Range("G11").Value = 3.825779881708576
$ws.Range("I11").Value = 44.65843329310331
$ws.Range("J11").Value = 7.465437757647796
$ws.Range("K11").Value = 22.71128885037693
$ws.Range("L11").Value = 14.7950284409154
$ws.Range("N11").Value = 28.06871754783284
$ws.Range("B12").Value = 25.26272016963065
$ws.Range("C12").Value = 11.26125481608721
$ws.Range("D12").Value = 12.87523007913676
$ws.Range("E12").Value = 11.97203397629409
$ws.Range("G12").Value = 3.82486614784374
$ws.Range("I12").Value = 44.67670801567706
$ws.Range("J12").Value = 7.464638663829557
$ws.Range("K12").Value = 22.73795285699373
$ws.Range("L12").Value = 14.79516077713735
$ws.Range("N12").Value = 28.06715656699436
$ws.Range("B13").Value = 25.254912090118
$ws.Range("C13").Value = 11.25328392453548
$ws.Range("D13").Value = 12.87440720060237
$ws.Range("E13").Value = 11.97288488586462
$ws.Range("G13").Value = 3.825062181973065
$ws.Range("I13").Value = 44.67273336115499
$ws.Range("J13").Value = 7.464809827507432
$ws.Range("K13").Value = 22.73217931792685
$ws.Range("L13").Value = 14.79511700886328
$ws.Range("N13").Value = 28.06748094105611
$ws.Range("B14").Value = 25.22958183248794
$ws.Range("C14").Value = 11.22732230994624
$ws.Range("D14").Value = 12.8717809066507
$ws.Range("E14").Value = 11.97569514427349
$ws.Range("G14").Value = 3.825704367994073
$ws.Range("I14").Value = 44.65991924463192
$ws.Range("J14").Value = 7.465371594061651
$ws.Range("K14").Value = 22.71346853768446
$ws.Range("L14").Value = 14.79503244751791
$ws.Range("N14").Value = 28.06858380273641
$ws.Range("B15").Value = 25.2141693761711
$ws.Range("C15").Value = 11.21144614286592
$ws.Range("D15").Value = 12.87021628133967
$ws.Range("E15").Value = 11.97744360783567
$ws.Range("G15").Value = 3.826099937387175
$ws.Range("I15").Value = 44.6521841155685
$ws.Range("J15").Value = 7.465718432807375
$ws.Range("K15").Value = 22.70209857664435
$ws.Range("L15").Value = 14.79502536773914
$ws.Range("N15").Value = 28.06929392116587
$ws.Range("B16").Value = 25.12744298424696
$ws.Range("C16").Value = 11.12088912225016
$ws.Range("D16").Value = 12.86192655176078
$ws.Range("E16").Value = 11.98787897701791
$ws.Range("G16").Value = 3.828400764797391
$ws.Range("I16").Value = 44.60960631374543
$ws.Range("J16").Value = 7.467748088219163
$ws.Range("K16").Value = 22.63834660995835
$ws.Range("L16").Value = 14.79567067154305
$ws.Range("N16").Value = 28.07389336932181
$ws.Range("B17").Value = 25.07568106722891
$ws.Range("C17").Value = 11.0657379445203
$ws.Range("D17").Value = 12.85744602392975
$ws.Range("E17").Value = 11.99465190864874
$ws.Range("G17").Value = 3.829842651429895
$ws.Range("I17").Value = 44.58505060473852
$ws.Range("J17").Value = 7.469030813021595
$ws.Range("K17").Value = 22.60050185916383
$ws.Range("L17").Value = 14.79667905986853
$ws.Range("N17").Value = 28.07718791293268
$ws.Range("B18").Value = 25.04644570463712
$ws.Range("C18").Value = 11.0341744143068
$ws.Range("D18").Value = 12.85509174558966
$ws.Range("E18").Value = 11.9986840391251
$ws.Range("G18").Value = 3.830683192780699
$ws.Range("I18").Value = 44.57150309058622
$ws.Range("J18").Value = 7.469782444471716
$ws.Range("K18").Value = 22.57920391404391
$ws.Range("L18").Value = 14.79748433587443
$ws.Range("N18").Value = 28.07925662952205
$ws.Range("B19").Value = 25.03663992713952
$ws.Range("C19").Value = 11.02351562316691
$ws.Range("D19").Value = 12.85433293724943
$ws.Range("E19").Value = 12.00007270732367
$ws.Range("G19").Value = 3.830969713474114
$ws.Range("I19").Value = 44.5670152347362
$ws.Range("J19").Value = 7.470039313765943
$ws.Range("K19").Value = 22.57207381892511
$ws.Range("L19").Value = 14.79779569569973
$ws.Range("N19").Value = 28.07998690660998
$ws.Range("B20").Value = 25.08113581411082
$ws.Range("C20").Value = 11.07159275765184
$ws.Range("D20").Value = 12.85789993499971
$ws.Range("E20").Value = 11.99391679186957
$ws.Range("G20").Value = 3.829688001077984
$ws.Range("I20").Value = 44.5876049890441
$ws.Range("J20").Value = 7.468892832747274
$ws.Range("K20").Value = 22.60448201972129
$ws.Range("L20").Value = 14.7965484036588
$ws.Range("N20").Value = 28.07681921701638
$ws.Range("B21").Value = 25.23700050614875
$ws.Range("C21").Value = 11.23494248790233
$ws.Range("D21").Value = 12.87254313479731
$ws.Range("E21").Value = 11.97486405648464
$ws.Range("G21").Value = 3.825515281615301
$ws.Range("I21").Value = 44.66365933588794
$ws.Range("J21").Value = 7.465206018626482
$ws.Range("K21").Value = 22.71894542380789
$ws.Range("L21").Value = 14.79504796759061
$ws.Range("N21").Value = 28.06825265822466
$ws.Range("B22").Value = 25.3434983719265
$ws.Range("C22").Value = 11.34289486454943
$ws.Range("D22").Value = 12.88408648640748
$ws.Range("E22").Value = 11.96362519881079
$ws.Range("G22").Value = 3.822887257135587
$ws.Range("I22").Value = 44.71846638574692
$ws.Range("J22").Value = 7.462919199126761
$ws.Range("K22").Value = 22.79783602241169
$ws.Range("L22").Value = 14.79606926204063
$ws.Range("N22").Value = 28.06420176172126
$ws.Range("B23").Value = 25.28624212436983
$ws.Range("C23").Value = 11.28517971785843
$ws.Range("D23").Value = 12.87774560852521
$ws.Range("E23").Value = 11.96951270297745
$ws.Range("G23").Value = 3.824280850490491
$ws.Range("I23").Value = 44.68874967548829
$ws.Range("J23").Value = 7.464128513944751
$ws.Range("K23").Value = 22.75536202790886
$ws.Range("L23").Value = 14.79534123026266
$ws.Range("N23").Value = 28.0662221619592
$ws.Range("B24").Value = 25.07866809154796
$ws.Range("C24").Value = 11.06894534693174
$ws.Range("D24").Value = 12.8576940313351
$ws.Range("E24").Value = 11.9942487074144
$ws.Range("G24").Value = 3.829757882418411
$ws.Range("I24").Value = 44.58644837635416
$ws.Range("J24").Value = 7.468955169469091
$ws.Range("K24").Value = 22.60268115676069
$ws.Range("L24").Value = 14.79660677076277
$ws.Range("N24").Value = 28.07698536040724
$ws.Range("B25").Value = 24.87535883451672
$ws.Range("C25").Value = 10.84186392689565
$ws.Range("D25").Value = 12.84458725724699
$ws.Range("E25").Value = 12.02615124436332
$ws.Range("G25").Value = 3.836094622698817
$ws.Range("I25").Value = 44.49814358442902
$ws.Range("J25").Value = 7.474691045831312
$ws.Range("K25").Value = 22.45598061145866
$ws.Range("L25").Value = 14.80658567199585
$ws.Range("N25").Value = 28.09523873938116
